# Application: Download: Add SpeedyPage
# Appends six new rows (148-153) of download-test data to Sheet1,
# matching the rows already recorded for the other routes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 148 --------------------------------------------------------------
$ws.Cells.Item(148, 1).Value = 45728
$ws.Cells.Item(148, 2).Value = 0.82500000000000007
$ws.Cells.Item(148, 3).Value = "阿里云北京H - Tokyo - SpeedyPage"
$ws.Cells.Item(148, 4).Value = "wget"
$ws.Cells.Item(148, 6).Value = 0.005

# --- Row 149 ----------------------------------------------------------------
$ws.Cells.Item(149, 2).Value = 0.82986111111111116
$ws.Cells.Item(149, 3).Value = "Azure JP - Tokyo - SpeedyPage"
$ws.Cells.Item(149, 4).Value = "wget"
$ws.Cells.Item(149, 6).Value = 66.2
$ws.Cells.Item(149, 9).Value = "2.7ms"

# --- Row 150 ----------------------------------------------------------------
$ws.Cells.Item(150, 2).Value = 0.83819444444444446
$ws.Cells.Item(150, 3).Value = "China - Azure JP - Tokyo - SpeedPage"
$ws.Cells.Item(150, 9).Value = "130.7ms"

# --- Row 151 ----------------------------------------------------------------
$ws.Cells.Item(151, 2).Value = 0.85277777777777775
$ws.Cells.Item(151, 3).Value = "阿里云北京H - Singapore - SpeedyPage"
$ws.Cells.Item(151, 4).Value = "wget"
$ws.Cells.Item(151, 6).Value = 5.94
$ws.Cells.Item(151, 9).Value = "125ms"

# --- Row 152 ----------------------------------------------------------------
$ws.Cells.Item(152, 2).Value = 0.85416666666666663
$ws.Cells.Item(152, 3).Value = "阿里云北京H - Los Angeles - SpeedyPage"
$ws.Cells.Item(152, 4).Value = "wget"
$ws.Cells.Item(152, 6).Value = 1.7
$ws.Cells.Item(152, 9).Value = "280ms"

# --- Row 153 ----------------------------------------------------------------
$ws.Cells.Item(153, 2).Value = 0.85555555555555562
$ws.Cells.Item(153, 3).Value = "阿里云北京H - Ashburn - SpeedyPage"
$ws.Cells.Item(153, 4).Value = "wget"
$ws.Cells.Item(153, 6).Value = 1.26
$ws.Cells.Item(153, 9).Value = "275ms"

# --- Formatting: reuse the existing Date/Time number formats --------------
# (column F already carries its number format at the column level, so only
#  columns A and B need their style copied explicitly, matching the style
#  used by the rows directly above.)
$ws.Range("A143").Copy()
$ws.Cells.Item(148, 1).PasteSpecial(-4122)

$ws.Range("B147").Copy()
$ws.Range("B148:B153").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View state: keep the selection / scroll position in sync with the
#     freshly-added rows (same relative offset as before the edit). -------
$ws.Range("D156").Select()
$excel.ActiveWindow.ScrollRow = 136
$excel.ActiveWindow.ScrollColumn = 1
